$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear old block (rows 3-13) that will be restructured ---
$ws.Range("A3:H13").Clear()

# --- Row 2 header (H2 text unchanged) ---
$ws.Range("H2").Value = "Duración"

# Row 3
$ws.Range("A3").Value = "SI"
$ws.Range("B3").Value = "Corregir vbles rotas"
$ws.Range("C3").Value = "Dataset original"
$ws.Range("D3").Value = 906
$ws.Range("E3").Value = "datasets/competencia3_2022.csv.gz"
$ws.Range("F3").Value = "exp/EC_CA9060/dataset.csv.gz"
$ws.Range("G3").Value = "Machine Learing"
$ws.Range("H3").Value = "2 minutos"

# Row 4
$ws.Range("A4").Value = "SI"
$ws.Range("B4").Value = "Corregir drifting"
$ws.Range("D4").Value = 914
$ws.Range("E4").Value = "exp/EC_CA9060/dataset.csv.gz"
$ws.Range("F4").Value = "exp/EC_DR9141/dataset.csv.gz"
$ws.Range("G4").Value = "Rank cero fijo."
$ws.Range("H4").Value = "10 minutos"

# Row 5
$ws.Range("A5").Value = "SI"
$ws.Range("B5").Value = "FE histórico y canarios"
$ws.Range("D5").Value = 925
$ws.Range("E5").Value = "exp/EC_DR9141/dataset.csv.gz"
$ws.Range("F5").Value = "exp/EC_FE9251/dataset.csv.gz"
$ws.Range("G5").Value = "Canarios asesinos. 257 variables finales"
$ws.Range("H5").Value = "15 minutos"

# Row 9
$ws.Range("B9").Value = "ALTERNATIVA 1: 2 meses"

# Row 10
$ws.Range("B10").Value = "Training strategy"
$ws.Range("D10").Value = "914_1"
$ws.Range("E10").Value = "exp/EC_FE9251/dataset.csv.gz"
$ws.Range("F10").Value = "exp/EC_TS931_1/dataset.csv.gz"
$ws.Range("G10").Value = "train & final_train: jun-19 a dic-19 y sep-20 a ene-21`nvalidation: feb-21`ntest: may-19 y may-21`nseed: 335897`nundersampling: 0.4"

# Row 11
$ws.Range("B11").Value = "BO"
$ws.Range("D11").Value = "942_1"
$ws.Range("E11").Value = "exp/EC_TS931_1/dataset.csv.gz"
$ws.Range("F11").Value = "exp/EC_HT9420_1/dataset.csv.gz"
$ws.Range("G11").Value = "Undersampling aplicado.`nKBO_iteraciones: 100`nseed: 335897"

# Row 15
$ws.Range("B15").Value = "ALTERNATIVA 2: 3 meses"

# Row 16
$ws.Range("B16").Value = "Rank para tratar Data Drifting"
$ws.Range("D16").Value = "914_2"
$ws.Range("E16").Value = "exp/EC_FE9251/dataset.csv.gz"
$ws.Range("F16").Value = "exp/EC_TS931_2/dataset.csv.gz"
$ws.Range("G16").Value = "train & final_train: may-19 a dic-19 y sep-20 a ene-21`nvalidation: feb-21`ntest: mar-21, abr-21, may-21.`nseed: 335897`nundersampling: 0.4"

# Row 17
$ws.Range("B17").Value = "BO"
$ws.Range("D17").Value = "942_2"
$ws.Range("E17").Value = "exp/EC_TS931_2/dataset.csv.gz"
$ws.Range("F17").Value = "exp/EC_HT9420_2/dataset.csv.gz"
$ws.Range("G17").Value = "Undersampling aplicado.`nKBO_iteraciones: 100`nseed: 335897"

# --- Styling: bold section headers (same bold style already used by former B7/B11) ---
$ws.Range("B9").Font.Bold = $true
$ws.Range("B15").Font.Bold = $true

# --- Styling: wrap text for the long multi-line observation cells ---
$ws.Range("G10").WrapText = $true
$ws.Range("G11").WrapText = $true
$ws.Range("G16").WrapText = $true
$ws.Range("G17").WrapText = $true

# --- Row heights for the wrapped rows ---
$ws.Rows.Item(10).RowHeight = 75
$ws.Rows.Item(11).RowHeight = 45
$ws.Rows.Item(16).RowHeight = 75
$ws.Rows.Item(17).RowHeight = 45

# --- Column widths (F, G widened; H newly added) ---
$ws.Columns.Item(6).ColumnWidth = 29.666666666666668
$ws.Columns.Item(7).ColumnWidth = 47.666666666666664
$ws.Columns.Item(8).ColumnWidth = 9.833333333333334

# --- Selection moves to A5 ---
$ws.Range("A5").Select() | Out-Null
